# Applies the commit's edits to the K-Pop article:
#  1. Merge the split title run into a single run (Find/Replace over the
#     whole title text forces Word to collapse it into one run).
#  2. A handful of small wording tweaks scattered through the body text.
#  3. Relocate the "_GoBack" bookmark from the title paragraph to the start
#     of the right-aligned "Autor: Martyna Dyba" paragraph.
#  4. Tweak line-spacing on the "Normal" and "Intense Quote" styles.

$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $ok = $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, `
                                   $true, 1, $false, $replace, 2)
    if (-not $ok) {
        Write-Output "WARNING: replace failed for: $find"
    }
}

# --- Title: collapse the three runs back into a single run -----------------
Replace-Text "K-Pop: Magia Rytmów i Kolorów, czyli Wprowadzenie do Koreańskiego Popu!" `
             "K-Pop: Magia Rytmów i Kolorów, czyli Wprowadzenie do Koreańskiego Popu!"

# --- Wording tweaks ----------------------------------------------------------
Replace-Text "czyli K-Popu! Jeśli jesteś ciekawym" "czyli k-popu! Jeśli jesteś ciekawym"

Replace-Text "K-Pop to krótka nazwa" "K-Pop to skrócona nazwa"

Replace-Text "Ikony K-Popu i Ich Showbiznesowe Podboje" "Ikony K-Popu i Ich Show biznesowe Podboje"

Replace-Text "koreańskim popem czy dopiero zaczynasz swoją przygodę – K-Pop to prawdziwa" `
             "koreańskim popem, czy dopiero zaczynasz swoją przygodę z k-popem to prawdziwa"

Replace-Text "Ciekawe jest to, że K-Pop nie zna granic!" "Ciekawe jest to, że ten gatunek nie zna granic!"

Replace-Text "To prawdziwie artystyczna miszmaszanka!" "To prawdziwie artystyczna mieszanka!"

Replace-Text "trendy w modzie K-Popowej, a nawet" "trendy w modzie koreańskiej, a nawet"

# --- Move the _GoBack bookmark to the start of the "Autor" paragraph -------
$authorPara = $d.Paragraphs.Last
$startRange = $d.Range($authorPara.Range.Start, $authorPara.Range.Start)
$d.Bookmarks.Add("_GoBack", $startRange)

# --- Style line-spacing tweaks -----------------------------------------------
$normal = $d.Styles("Normal")
$normal.ParagraphFormat.LineSpacingRule = 5
$normal.ParagraphFormat.LineSpacing = 12.8

$quote = $d.Styles("Intense Quote")
$quote.ParagraphFormat.LineSpacingRule = 5
$quote.ParagraphFormat.LineSpacing = 12.95

Write-Output "done"
